$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new blank rows before row 192, pushing existing rows 192:279 down to 195:282.
$ws.Rows("192:194").Insert()

# Row 192 - new data entry
$ws.Range("A192").Value = 10
$ws.Range("B192").Value = "Vega Modelo de Temuco"
$ws.Range("C192").Value = "La Araucanía"
$ws.Range("D192").Value = 44875
$ws.Range("E192").Value = 9
$ws.Range("F192").Value = "Fruta"
$ws.Range("G192").Value = 100101
$ws.Range("H192").Value = "Berries"
$ws.Range("I192").Value = 100112025
$ws.Range("J192").Value = "Frutilla"
$ws.Range("K192").Value = "Sin especificar"
$ws.Range("L192").Value = "Primera"
$ws.Range("M192").Value = 5000
$ws.Range("N192").Value = 8000
$ws.Range("O192").Value = 9000
$ws.Range("P192").Value = 8400
$ws.Range("Q192").Value = '$/bandeja 7 kilos'
$ws.Range("R192").Value = "Provincia de Melipilla"
$ws.Range("S192").Value = 1200
$ws.Range("T192").Value = 7

# Row 193 - new data entry
$ws.Range("A193").Value = 10
$ws.Range("B193").Value = "Vega Modelo de Temuco"
$ws.Range("C193").Value = "La Araucanía"
$ws.Range("D193").Value = 44875
$ws.Range("E193").Value = 9
$ws.Range("F193").Value = "Fruta"
$ws.Range("G193").Value = 100101
$ws.Range("H193").Value = "Berries"
$ws.Range("I193").Value = 100112025
$ws.Range("J193").Value = "Frutilla"
$ws.Range("K193").Value = "Sin especificar"
$ws.Range("L193").Value = "Primera"
$ws.Range("M193").Value = 500
$ws.Range("N193").Value = 8000
$ws.Range("O193").Value = 8000
$ws.Range("P193").Value = 8000
$ws.Range("Q193").Value = '$/caja 7 kilos'
$ws.Range("R193").Value = "Región de La Araucanía"
$ws.Range("S193").Value = 1143
$ws.Range("T193").Value = 7

# Row 194 - new data entry
$ws.Range("A194").Value = 10
$ws.Range("B194").Value = "Vega Modelo de Temuco"
$ws.Range("C194").Value = "La Araucanía"
$ws.Range("D194").Value = 44875
$ws.Range("E194").Value = 9
$ws.Range("F194").Value = "Fruta"
$ws.Range("G194").Value = 100101
$ws.Range("H194").Value = "Berries"
$ws.Range("I194").Value = 100112025
$ws.Range("J194").Value = "Frutilla"
$ws.Range("K194").Value = "Sin especificar"
$ws.Range("L194").Value = "Segunda"
$ws.Range("M194").Value = 500
$ws.Range("N194").Value = 6000
$ws.Range("O194").Value = 6500
$ws.Range("P194").Value = 6300
$ws.Range("Q194").Value = '$/bandeja 7 kilos'
$ws.Range("R194").Value = "Provincia de Melipilla"
$ws.Range("S194").Value = 900
$ws.Range("T194").Value = 7
